$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

$s.Shapes.Item("TextBox 23").Delete()
$s.Shapes.Item("Straight Arrow Connector 22").Delete()
$s.Shapes.Item("Picture 2").Delete()
$s.Shapes.Item("Picture 19").Delete()
